$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 0.06981456006938053, 7.819040571458529),
    @(3, 0.07418709268332681, 8.45262065420218),
    @(4, -0.03662410462721217, 18.41918759872105),
    @(5, -0.02827920324327398, -4.459839464958129),
    @(6, -0.09375036274816409, 11.59062305607925),
    @(7, -0.0905156173356471, 0.9363463242187672),
    @(8, -0.3522130331728875, 4.004174130422792),
    @(9, -0.3707903866088169, 4.947653685622434),
    @(10, 0.02243065816861043, 11.09237022080871),
    @(11, 0.04045800286561223, 78.27246941048899),
    @(12, 0.2171903122419736, -2.055147260158437),
    @(13, 0.2286198763218104, 1.515529248309028),
    @(14, -0.0423543248488569, -0.5919048331448972),
    @(15, -0.0453565555941181, 4.905498039464362),
    @(16, 0.212957049951566, 0.1804793927103588),
    @(17, 0.2195398528038392, -0.4570940385320715),
    @(18, 0.08065590735726356, 10.45691811470875),
    @(19, 0.07731856556847012, 2.628889064839274),
    @(20, -0.08130991902103457, -8.439285331583291),
    @(21, -0.08667399836390689, -0.1214836521738738),
    @(22, 0.06117531173446256, -16.77043128219126),
    @(23, 0.07547849250595139, 10.45959090777678),
    @(24, 0.05919867730843299, -11.12928493580626),
    @(25, 0.07076908915051398, 29.18769525110326),
    @(26, 0.1188644250137464, -0.4060436773159631),
    @(27, 0.1147161213579321, 0.7517448496137927),
    @(28, 0.1334537448415444, 3.24971461008987),
    @(29, 0.1493048118165373, -1.018030426905561),
    @(30, 0.08660190648115718, 2.721708996884293),
    @(31, 0.08422565800560795, 3.106468377266105),
    @(32, 0.0578300106329763, 8.377290403527534),
    @(33, 0.05327228878183658, -3.567194456656766),
    @(34, 0.01685976678628124, -2.866107803355784),
    @(35, 0.01476290257631001, -12.64910365074669),
    @(36, -0.02658450229724496, 8.473537083859583),
    @(37, -0.02339515210720308, 29.66799711960084),
    @(38, 0.07651656826372312, -2.256123270025682),
    @(39, 0.07566226998961434, -2.679088597869588),
    @(40, 0.06778841244179645, 2.385446177579617),
    @(41, 0.07014672341214186, 7.880129505104935),
    @(42, 0.08097719768867399, 4.095451345300171),
    @(43, 0.08002140611715626, -0.1784904084086701),
    @(44, 0.07881683934424609, -10.6864034650987),
    @(45, 0.08983344264494685, -0.6085877808285041),
    @(46, 0.01192697999351138, 535.8923534729275),
    @(47, -0.0006092236977627549, -1175.236313981536),
    @(48, -0.09717906322142596, -1.117441703298478),
    @(49, -0.1007602460744671, 8.045503291318386),
    @(50, 0.1710738953855934, 0.3370128035295785),
    @(51, 0.1739833048481258, 2.442629277708561),
    @(52, 0.06452388715755088, -9.072472456051909),
    @(53, 0.06618202438924688, 2.901164808447063),
    @(54, -0.1318780570844196, -3.188662266730049),
    @(55, -0.1342217989502003, -15.24036525875476),
    @(56, 0.1893462526135997, -0.3550044977942086),
    @(57, 0.195815146103937, -1.555052582233344),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 7).Value = $row[1]
    $ws.Cells.Item($r, 8).Value = $row[2]
}

$ws.Range("I2").Value = -8.241209273008328
